# Applies:
#  1) slide5.xml: the table's style id changes from
#     {1F9C4039-A4B1-4A02-97CF-F9B796A9239E} to
#     {40EDE20E-B1D6-4E89-B748-BFE56740679E}
#  2) the deck's theme colour scheme (the one backing ppt/theme/theme2.xml,
#     which is the theme actually wired to the slide master / presentation)
#     is swapped from the "Integral / Red Violet" palette to the
#     "Office Theme / Office" palette that used to live in ppt/theme/theme1.xml.

$p = $ppt.ActivePresentation

# --- 1. Table style id (slide 5, shape 2 - the graphicFrame/table) ---------
$s = $p.Slides.Item(5)
$tblShape = $s.Shapes.Item(2)
$tblShape.Table.ApplyStyle("{40EDE20E-B1D6-4E89-B748-BFE56740679E}")

# --- 2. Swap the active theme's 12 colours to the old "Office Theme" ------
# Order exposed by ThemeColorScheme.Colors(1..12):
#   dk1, lt1, dk2, lt2, accent1..accent6, hlink, folHlink
$officeThemeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x6A5444,  # dk2      (44546A)
    0xE6E6E7,  # lt2      (E7E6E6)
    0xD59B5B,  # accent1  (5B9BD5)
    0x317DED,  # accent2  (ED7D31)
    0xA5A5A5,  # accent3  (A5A5A5)
    0x00C0FF,  # accent4  (FFC000)
    0xC47244,  # accent5  (4472C4)
    0x47AD70,  # accent6  (70AD47)
    0xC16305,  # hlink    (0563C1)
    0x724F95   # folHlink (954F72)
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeColors[$i - 1]
}
